# Update the "today's snapshot" row (row 2) with refreshed quote data.
# Force each target cell to Text format BEFORE writing so Excel does not
# silently reinterpret numeric-looking strings (e.g. "28.9") as numbers -
# the source file stores every data cell as literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2Cells = @("B2", "C2", "D2", "F2", "G2")
foreach ($addr in $row2Cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("B2").Value = "$71.43"
$ws.Range("C2").Value = "+0.67(0.95%) 1D"
$ws.Range("D2").Value = "$307.22B"
$ws.Range("F2").Value = "28.9"
$ws.Range("G2").Value = "12.4"

foreach ($addr in $row2Cells) {
    $ws.Range($addr).Style = "Normal"
}

# Add the new "Cash Flow" block: a bold/bordered header row at row 10
# (mirroring the other section headers in rows 1/4/7) followed by a data
# row at row 11. Row 9 is left blank, matching the existing separator-row
# pattern between sections.
$ws.Range("A10").Value = "Investing Cash Flow"
$ws.Range("B10").Value = "Operations Cash Flow"
$ws.Range("C10").Value = "Financing Cash Flow"
$ws.Range("D10").Value = "Net Cash Flow"
$ws.Range("E10").Value = "Free Cash Flow"
$ws.Range("F10").Value = "Capital Expenditure"
$ws.Range("G10").Value = "Cash and Equivalents"
$ws.Range("H10").Value = "Payments & Cash Distribution"
$ws.Range("I10").Value = "Basic Common Share"
$ws.Range("J10").Value = "Working Capital"

# Copy the formatting (bold font, borders, centered alignment) from the
# existing header row 1 onto the new header row 10.
$ws.Range("A1:J1").Copy()
$ws.Range("A10:J10").PasteSpecial(-4122)

# Row 11: the cash-flow figures, entered as literal text (forcing Text
# format first) so values like "-763" / "11,018" keep their comma
# formatting instead of becoming numbers.
$dataRow11 = $ws.Range("A11:J11")
$dataRow11.NumberFormat = "@"

$ws.Range("A11").Value = "-763"
$ws.Range("B11").Value = "11,018"
$ws.Range("C11").Value = "-10,250"
$ws.Range("D11").Value = "-200"
$ws.Range("E11").Value = "9,609"
$ws.Range("F11").Value = "-1,409"
$ws.Range("G11").Value = "10,562"
$ws.Range("H11").Value = "-7,616"
$ws.Range("I11").Value = "0"
$ws.Range("J11").Value = "2,867"

$dataRow11.Style = "Normal"
